$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Étlap rögzítése")

# Update the text of the two requirement cells (C12, C13)
$ws.Range("C12").Value = "A rendszer megjeleníti az étlap rögzítésénél beállítható paramétereket: kategória neve, étel neve,  ára, akciós ára, akció időszaka"
$ws.Range("C13").Value = "Az étteremvezető kiválasztja a megfelelő paramétereket: kategória neve, étel neve,  ára, akciós ára, akció időszaka"

# Adjust row heights for the two rows to accommodate the longer text
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 45

# Update the view: scroll so row 6 is the top-left visible row, and move the
# active selection to F13
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("F13").Select()
